# Apply the commit: insert a comma after PROPER(Bn) in the CONCATENATE
# formula of column H (rows 2-59) on the "Datos" worksheet, and move the
# active cell selection from H63 to H10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

for ($r = 2; $r -le 59; $r++) {
    $formula = '=CONCATENATE("Yo, ",PROPER(B' + $r + '),", con DNI ",A' + $r + ',' + `
        '" y con el grupo profesional ",VLOOKUP(E' + $r + ',grupos,2,FALSE),' + `
        '" percibo un salario de ",Datos!F' + $r + ',' + '" pesetas.")'
    $ws.Range("H$r").Formula = $formula
}

# Update the selection shown in the sheet view to H10
$ws.Activate()
$ws.Range("H10").Select()
